$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 9 ("Golang Development Frameworks"): the "GXUI, Fyne" paragraph
# had a stray trailing <a:endParaRPr> left over from an earlier edit;
# drop it by deleting that paragraph and re-inserting identical text as
# a fresh paragraph (fresh paragraphs come back with no endParaRPr).
# ---------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$tr9 = $slide9.Shapes.Item(2).TextFrame.TextRange

$gxuiPara = $tr9.Paragraphs(2, 1)
$gxuiPara.Delete()

$goMobilePara = $tr9.Paragraphs(2, 1)
$goMobilePara.InsertBefore("GXUI, Fyne" + [char]13)

# ---------------------------------------------------------------------
# Slide 10 ("Up next:"): extend the "Create: GUI app using Fyne" line
# with " Or Qt", and add a brand-new "Create: Mobile App using GoMobile"
# bullet right before the trailing blank paragraph.
# ---------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$tr10 = $slide10.Shapes.Item(2).TextFrame.TextRange

$guiPara = $tr10.Paragraphs(1, 1)
$guiPara.Text = "Create: GUI app using Fyne Or Qt"
$guiPara.Characters(1, 22).Text = "Create: GUI app using "
$guiPara.Characters(23, 8).Text = "Fyne Or "
$guiPara.Characters(31, 2).Text = "Qt"

$lastPara = $tr10.Paragraphs(4, 1)
$lastPara.InsertBefore("Create: Mobile Ap" + "p using " + "GoMobile" + [char]13)
